$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the lower block of rows (old 13-17) down by 6 rows so it lands on 19-23,
# leaving rows 12-18 empty (mirrors the target layout).
$ws.Range("B13:D18").EntireRow.Insert()

# The old blank spacer row (previously row 11, still styled) now sits at row 11;
# strip its leftover formatting since it will be replaced with fresh content below.
$ws.Range("B11:D11").ClearFormats()

# Row 9: new feature entry, styled like the existing "Garrett" rows (strikethrough font, style index 2)
$ws.Range("B9:D9").Font.Strikethrough = $true
$ws.Range("B9").Value = "Ball serve location based on serve quality"
$ws.Range("C9").Value = "Garrett"
$ws.Range("D9").Value = "When the player or AI serves the ball, if the quality is low the ball position can be changed randomly"

# Row 10: new feature entry, plain font (no strikethrough)
$ws.Range("B10").Value = "AI chooses where to attack"
$ws.Range("C10").Value = "Garrett"
$ws.Range("D10").Value = "Have the AI choose randomly where to attack, to be improved in future"

# Row 11: new feature entry, plain font (no strikethrough)
$ws.Range("B11").Value = "Ball attack location change based on quality"
$ws.Range("C11").Value = "Garrett"
$ws.Range("D11").Value = "The location of the ball can change randomly based on the quality of the attack, if the attack randomly goes off the grid, it is an error. Need to also reduce the chances for hitting errors to account for these new error chances"

$ws.Range("D12").Select()
